$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "Adam Milne"
$ws.Name = "Adam Milne"

# Final 4x13 grid (header + 3 data rows). A new leading "matchNo" column was
# inserted, shifting every previously-existing column one place to the
# right, and two additional match rows were scraped in.
$grid = @(
    @("matchNo","teamName","batterName","states","runs","balls","fours","sixes","sr","opponentTeamName","venue","date","result"),
    @("30th","Mumbai Indians","Adam Milne","c sub (K Gowtham) b Bravo","15","15","0","1","100.00","Chennai Super Kings","Dubai (DSC)","September 19","Super Kings won by 20 runs"),
    @("39th","Mumbai Indians","Adam Milne","b Patel","0","1","0","0","0.00","Royal Challengers Bangalore","Dubai (DSC)","September 26","RCB won by 54 runs"),
    @("34th","Mumbai Indians","Adam Milne","","1","1","0","0","100.00","Kolkata Knight Riders","Abu Dhabi","September 23","KKR won by 7 wickets (with 29 balls remaining)")
)

for ($r = 0; $r -lt $grid.Length; $r++) {
    $rowData = $grid[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $cell = $ws.Cells.Item($r + 1, $c + 1)
        # Force text storage (matches the source file's t="str" cells) so
        # values like "100.00" or "0" aren't auto-coerced into numbers.
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[$c]
    }
}
